# Update the Wnt2-Fzd4 NATMI LR-pair sheet with recomputed TPM values and
# two new target clusters (Inflammatory-Mac, Resolving-Mac).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2 (Sending: FAPs -> Target: ECs) - recomputed M:T values
# ---------------------------------------------------------------------
$ws.Range("M2").Value = 29.75868033333333
$ws.Range("N2").Value = 89.27604099999999
$ws.Range("O2").Value = 0.4948552779010537
$ws.Range("P2").Value = 0.4948552779010535
$ws.Range("Q2").Value = 1.637868167746111
$ws.Range("R2").Value = 14.740813509715
$ws.Range("S2").Value = 0.4948552779010537
$ws.Range("T2").Value = 0.4948552779010535

# ---------------------------------------------------------------------
# Row 3 (Sending: FAPs -> Target: FAPs) - recomputed O,P,S,T values
# ---------------------------------------------------------------------
$ws.Range("O3").Value = 0.2919251856942525
$ws.Range("P3").Value = 0.2919251856942524
$ws.Range("S3").Value = 0.2919251856942525
$ws.Range("T3").Value = 0.2919251856942524

# ---------------------------------------------------------------------
# Row 4 (Sending: FAPs -> Target: Inflammatory-Mac, new cluster taking
# this row's former "MuSCs" slot) - full K:T rewrite
# ---------------------------------------------------------------------
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1421396666666667
$ws.Range("N4").Value = 0.426419
$ws.Range("O4").Value = 0.002363631836533717
$ws.Range("P4").Value = 0.002363631836533717
$ws.Range("Q4").Value = 0.00782313035388889
$ws.Range("R4").Value = 0.07040817318500001
$ws.Range("S4").Value = 0.002363631836533717
$ws.Range("T4").Value = 0.002363631836533717

# ---------------------------------------------------------------------
# Row 5 (Sending: FAPs -> Target: MuSCs) - new row
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt2"
$ws.Range("C5").Value = "Fzd4"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.05503833333333333
$ws.Range("H5").Value = 0.165115
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 12.42872866666667
$ws.Range("N5").Value = 37.286186
$ws.Range("O5").Value = 0.2066765699758167
$ws.Range("P5").Value = 0.2066765699758166
$ws.Range("Q5").Value = 0.6840565112655556
$ws.Range("R5").Value = 6.156508601390001
$ws.Range("S5").Value = 0.2066765699758167
$ws.Range("T5").Value = 0.2066765699758166

# ---------------------------------------------------------------------
# Row 6 (Sending: FAPs -> Target: Resolving-Mac) - new row
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt2"
$ws.Range("C6").Value = "Fzd4"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.05503833333333333
$ws.Range("H6").Value = 0.165115
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.251329
$ws.Range("N6").Value = 0.753987
$ws.Range("O6").Value = 0.004179334592343558
$ws.Range("P6").Value = 0.004179334592343557
$ws.Range("Q6").Value = 0.01383272927833333
$ws.Range("R6").Value = 0.124494563505
$ws.Range("S6").Value = 0.004179334592343558
$ws.Range("T6").Value = 0.004179334592343557
